$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1729.5294  # H33: 1727.4117 -> 1729.5294
$ws.Cells.Item(33, 9).Value = 263.81818  # I33: 260.54544 -> 263.81818
$ws.Cells.Item(33, 11).Value = 263.81818  # K33: 260.54544 -> 263.81818
$ws.Cells.Item(33, 13).Value = -34.81817999999998  # M33: -31.54543999999999 -> -34.81817999999998
$ws.Cells.Item(74, 8).Value = 11059.2  # H74: 11276.889 -> 11059.2
$ws.Cells.Item(74, 10).Value = 14052  # J74: 19004 -> 14052
$ws.Cells.Item(74, 12).Value = 14052  # L74: 19004 -> 14052
$ws.Cells.Item(74, 14).Value = -15924  # N74: -20876 -> -15924
$ws.Cells.Item(76, 8).Value = 6600.25  # H76: 7949.25 -> 6600.25
$ws.Cells.Item(76, 9).Value = 3239.8  # I76: 4396.5 -> 3239.8
$ws.Cells.Item(76, 10).Value = 9000.571  # J76: 11502 -> 9000.571
$ws.Cells.Item(76, 11).Value = 3239.8  # K76: 4396.5 -> 3239.8
$ws.Cells.Item(76, 12).Value = 9000.571  # L76: 11502 -> 9000.571
$ws.Cells.Item(76, 13).Value = -2924.8  # M76: -4081.5 -> -2924.8
$ws.Cells.Item(76, 14).Value = -9630.571  # N76: -12132 -> -9630.571
$ws.Cells.Item(77, 8).Value = 11059.2  # H77: 11276.889 -> 11059.2
$ws.Cells.Item(77, 10).Value = 14052  # J77: 19004 -> 14052
$ws.Cells.Item(77, 12).Value = 70260  # L77: 95020 -> 70260
$ws.Cells.Item(77, 14).Value = -79620  # N77: -104380 -> -79620
$ws.Cells.Item(79, 8).Value = 6600.25  # H79: 7949.25 -> 6600.25
$ws.Cells.Item(79, 9).Value = 3239.8  # I79: 4396.5 -> 3239.8
$ws.Cells.Item(79, 10).Value = 9000.571  # J79: 11502 -> 9000.571
$ws.Cells.Item(79, 11).Value = 3239.8  # K79: 4396.5 -> 3239.8
$ws.Cells.Item(79, 12).Value = 9000.571  # L79: 11502 -> 9000.571
$ws.Cells.Item(79, 13).Value = -2147.8  # M79: -3304.5 -> -2147.8
$ws.Cells.Item(79, 14).Value = -11184.571  # N79: -13686 -> -11184.571
$ws.Cells.Item(86, 8).Value = 4977.6665  # H86: 4033 -> 4977.6665
$ws.Cells.Item(86, 9).Value = 3500  # I86: 4000 -> 3500
$ws.Cells.Item(86, 10).Value = 5399.857  # J86: 4049.5 -> 5399.857
$ws.Cells.Item(86, 11).Value = 3500  # K86: 4000 -> 3500
$ws.Cells.Item(86, 12).Value = 5399.857  # L86: 4049.5 -> 5399.857
$ws.Cells.Item(86, 13).Value = -2377  # M86: -2877 -> -2377
$ws.Cells.Item(86, 14).Value = -7645.857  # N86: -6295.5 -> -7645.857
$ws.Cells.Item(89, 8).Value = 4977.6665  # H89: 4033 -> 4977.6665
$ws.Cells.Item(89, 9).Value = 3500  # I89: 4000 -> 3500
$ws.Cells.Item(89, 10).Value = 5399.857  # J89: 4049.5 -> 5399.857
$ws.Cells.Item(89, 11).Value = 17500  # K89: 20000 -> 17500
$ws.Cells.Item(89, 12).Value = 26999.285  # L89: 20247.5 -> 26999.285
$ws.Cells.Item(89, 13).Value = -11884  # M89: -14384 -> -11884
$ws.Cells.Item(89, 14).Value = -38231.285  # N89: -31479.5 -> -38231.285
$ws.Cells.Item(98, 8).Value = 921.8  # H98: 927.25 -> 921.8
$ws.Cells.Item(98, 9).Value = 921.8  # I98: 927.25 -> 921.8
$ws.Cells.Item(98, 11).Value = 921.8  # K98: 927.25 -> 921.8
$ws.Cells.Item(98, 13).Value = 576.2  # M98: 570.75 -> 576.2
$ws.Cells.Item(103, 8).Value = 4741.9  # H103: 5477.375 -> 4741.9
$ws.Cells.Item(103, 9).Value = 1800  # I103: 0 -> 1800
$ws.Cells.Item(103, 11).Value = 5400  # K103: 0 -> 5400
$ws.Cells.Item(103, 13).Value = -4814  # M103: None -> -4814
$ws.Cells.Item(112, 8).Value = 1423.0435  # H112: 1337.9412 -> 1423.0435
$ws.Cells.Item(112, 9).Value = 1096.625  # I112: 1111.2858 -> 1096.625
$ws.Cells.Item(112, 10).Value = 1597.1333  # J112: 1496.6 -> 1597.1333
$ws.Cells.Item(112, 11).Value = 3289.875  # K112: 3333.8574 -> 3289.875
$ws.Cells.Item(112, 12).Value = 4791.3999  # L112: 4489.799999999999 -> 4791.3999
$ws.Cells.Item(112, 13).Value = -2181.875  # M112: -2225.8574 -> -2181.875
$ws.Cells.Item(112, 14).Value = -7007.3999  # N112: -6705.799999999999 -> -7007.3999
$ws.Cells.Item(122, 8).Value = 921.8  # H122: 927.25 -> 921.8
$ws.Cells.Item(122, 9).Value = 921.8  # I122: 927.25 -> 921.8
$ws.Cells.Item(122, 11).Value = 2765.4  # K122: 2781.75 -> 2765.4
$ws.Cells.Item(122, 13).Value = -315.3999999999996  # M122: -331.75 -> -315.3999999999996
$ws.Cells.Item(133, 8).Value = 61999.2  # H133: 62000 -> 61999.2
$ws.Cells.Item(133, 10).Value = 61999.2  # J133: 62000 -> 61999.2
$ws.Cells.Item(133, 12).Value = 61999.2  # L133: 62000 -> 61999.2
$ws.Cells.Item(133, 14).Value = -72119.2  # N133: -72120 -> -72119.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 34231.668  # H2: 12945.667 -> 34231.668
$ws.Cells.Item(2, 9).Value = 1347.5  # I2: 1347 -> 1347.5
$ws.Cells.Item(2, 10).Value = 100000  # J2: 18745 -> 100000
$ws.Cells.Item(2, 11).Value = 1347.5  # K2: 1347 -> 1347.5
$ws.Cells.Item(2, 12).Value = 100000  # L2: 18745 -> 100000
$ws.Cells.Item(2, 13).Value = -1234.5  # M2: -1234 -> -1234.5
$ws.Cells.Item(2, 14).Value = -100226  # N2: -18971 -> -100226
$ws.Cells.Item(45, 10).Value = 5200.6665  # J45: 5211.778 -> 5200.6665
$ws.Cells.Item(45, 12).Value = 5200.6665  # L45: 5211.778 -> 5200.6665
$ws.Cells.Item(45, 14).Value = -5954.6665  # N45: -5965.778 -> -5954.6665
$ws.Cells.Item(61, 8).Value = 6158.8057  # H61: 6156.0835 -> 6158.8057
$ws.Cells.Item(61, 9).Value = 4834.625  # I61: 4746.0605 -> 4834.625
$ws.Cells.Item(61, 10).Value = 16752.25  # J61: 21666.334 -> 16752.25
$ws.Cells.Item(61, 11).Value = 4834.625  # K61: 4746.0605 -> 4834.625
$ws.Cells.Item(61, 12).Value = 16752.25  # L61: 21666.334 -> 16752.25
$ws.Cells.Item(61, 13).Value = -4622.625  # M61: -4534.0605 -> -4622.625
$ws.Cells.Item(61, 14).Value = -17176.25  # N61: -22090.334 -> -17176.25
$ws.Cells.Item(63, 8).Value = 7325  # H63: 6755.5557 -> 7325
$ws.Cells.Item(63, 9).Value = 2866.6667  # I63: 2700 -> 2866.6667
$ws.Cells.Item(63, 11).Value = 2866.6667  # K63: 2700 -> 2866.6667
$ws.Cells.Item(63, 13).Value = -2180.6667  # M63: -2014 -> -2180.6667
$ws.Cells.Item(66, 8).Value = 7325  # H66: 6755.5557 -> 7325
$ws.Cells.Item(66, 9).Value = 2866.6667  # I66: 2700 -> 2866.6667
$ws.Cells.Item(66, 11).Value = 14333.3335  # K66: 13500 -> 14333.3335
$ws.Cells.Item(66, 13).Value = -10901.3335  # M66: -10068 -> -10901.3335
$ws.Cells.Item(74, 8).Value = 55559220  # H74: 47622764 -> 55559220
$ws.Cells.Item(74, 9).Value = 66670264  # I74: 55559224 -> 66670264
$ws.Cells.Item(74, 10).Value = 4000  # J74: 4014 -> 4000
$ws.Cells.Item(74, 11).Value = 66670264  # K74: 55559224 -> 66670264
$ws.Cells.Item(74, 12).Value = 4000  # L74: 4014 -> 4000
$ws.Cells.Item(74, 13).Value = -66669390  # M74: -55558350 -> -66669390
$ws.Cells.Item(74, 14).Value = -5748  # N74: -5762 -> -5748
$ws.Cells.Item(77, 8).Value = 55559220  # H77: 47622764 -> 55559220
$ws.Cells.Item(77, 9).Value = 66670264  # I77: 55559224 -> 66670264
$ws.Cells.Item(77, 10).Value = 4000  # J77: 4014 -> 4000
$ws.Cells.Item(77, 11).Value = 333351320  # K77: 277796120 -> 333351320
$ws.Cells.Item(77, 12).Value = 20000  # L77: 20070 -> 20000
$ws.Cells.Item(77, 13).Value = -333346952  # M77: -277791752 -> -333346952
$ws.Cells.Item(77, 14).Value = -28736  # N77: -28806 -> -28736
$ws.Cells.Item(97, 8).Value = 1109.8182  # H97: 1313.1111 -> 1109.8182
$ws.Cells.Item(97, 9).Value = 1210.7778  # I97: 1338.375 -> 1210.7778
$ws.Cells.Item(97, 10).Value = 655.5  # J97: 1111 -> 655.5
$ws.Cells.Item(97, 11).Value = 1210.7778  # K97: 1338.375 -> 1210.7778
$ws.Cells.Item(97, 12).Value = 655.5  # L97: 1111 -> 655.5
$ws.Cells.Item(97, 13).Value = -714.7778000000001  # M97: -842.375 -> -714.7778000000001
$ws.Cells.Item(97, 14).Value = -1647.5  # N97: -2103 -> -1647.5
$ws.Cells.Item(116, 8).Value = 34231.668  # H116: 12945.667 -> 34231.668
$ws.Cells.Item(116, 9).Value = 1347.5  # I116: 1347 -> 1347.5
$ws.Cells.Item(116, 10).Value = 100000  # J116: 18745 -> 100000
$ws.Cells.Item(116, 11).Value = 1347.5  # K116: 1347 -> 1347.5
$ws.Cells.Item(116, 12).Value = 100000  # L116: 18745 -> 100000
$ws.Cells.Item(116, 13).Value = 946.5  # M116: 947 -> 946.5
$ws.Cells.Item(116, 14).Value = -104588  # N116: -23333 -> -104588
$ws.Cells.Item(136, 8).Value = 6158.8057  # H136: 6156.0835 -> 6158.8057
$ws.Cells.Item(136, 9).Value = 4834.625  # I136: 4746.0605 -> 4834.625
$ws.Cells.Item(136, 10).Value = 16752.25  # J136: 21666.334 -> 16752.25
$ws.Cells.Item(136, 11).Value = 14503.875  # K136: 14238.1815 -> 14503.875
$ws.Cells.Item(136, 12).Value = 50256.75  # L136: 64999.00199999999 -> 50256.75
$ws.Cells.Item(136, 13).Value = -11953.875  # M136: -11688.1815 -> -11953.875
$ws.Cells.Item(136, 14).Value = -55356.75  # N136: -70099.00199999999 -> -55356.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 34231.668  # H3: 12945.667 -> 34231.668
$ws.Cells.Item(3, 9).Value = 1347.5  # I3: 1347 -> 1347.5
$ws.Cells.Item(3, 10).Value = 100000  # J3: 18745 -> 100000
$ws.Cells.Item(3, 11).Value = 1347.5  # K3: 1347 -> 1347.5
$ws.Cells.Item(3, 12).Value = 100000  # L3: 18745 -> 100000
$ws.Cells.Item(3, 13).Value = -1233.5  # M3: -1233 -> -1233.5
$ws.Cells.Item(3, 14).Value = -100228  # N3: -18973 -> -100228

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1343.4667  # H107: 902.9524 -> 1343.4667
$ws.Cells.Item(107, 9).Value = 1366.1666  # I107: 654.38464 -> 1366.1666
$ws.Cells.Item(107, 10).Value = 1328.3334  # J107: 1306.875 -> 1328.3334
$ws.Cells.Item(107, 11).Value = 1366.1666  # K107: 654.38464 -> 1366.1666
$ws.Cells.Item(107, 12).Value = 1328.3334  # L107: 1306.875 -> 1328.3334
$ws.Cells.Item(107, 13).Value = 553.8334  # M107: 1265.61536 -> 553.8334
$ws.Cells.Item(107, 14).Value = -5168.3334  # N107: -5146.875 -> -5168.3334
$ws.Cells.Item(134, 8).Value = 4161.2  # H134: 4701.875 -> 4161.2
$ws.Cells.Item(134, 9).Value = 2459.3333  # I134: 2689.75 -> 2459.3333
$ws.Cells.Item(134, 11).Value = 7377.999899999999  # K134: 8069.25 -> 7377.999899999999
$ws.Cells.Item(134, 13).Value = -4842.999899999999  # M134: -5534.25 -> -4842.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 15635.272  # H5: 16698.9 -> 15635.272
$ws.Cells.Item(5, 10).Value = 21717  # J5: 24503.334 -> 21717
$ws.Cells.Item(5, 12).Value = 65151  # L5: 73510.00199999999 -> 65151
$ws.Cells.Item(5, 14).Value = -65375  # N5: -73734.00199999999 -> -65375
$ws.Cells.Item(22, 8).Value = 8333.666999999999  # H22: 8501 -> 8333.666999999999
$ws.Cells.Item(22, 10).Value = 8333.666999999999  # J22: 8501 -> 8333.666999999999
$ws.Cells.Item(22, 12).Value = 25001.001  # L22: 25503 -> 25001.001
$ws.Cells.Item(22, 14).Value = -25339.001  # N22: -25841 -> -25339.001
$ws.Cells.Item(27, 8).Value = 8333.666999999999  # H27: 8501 -> 8333.666999999999
$ws.Cells.Item(27, 10).Value = 8333.666999999999  # J27: 8501 -> 8333.666999999999
$ws.Cells.Item(27, 12).Value = 25001.001  # L27: 25503 -> 25001.001
$ws.Cells.Item(27, 14).Value = -25205.001  # N27: -25707 -> -25205.001
$ws.Cells.Item(109, 8).Value = 1889.3334  # H109: 1934 -> 1889.3334
$ws.Cells.Item(135, 8).Value = 15635.272  # H135: 16698.9 -> 15635.272
$ws.Cells.Item(135, 10).Value = 21717  # J135: 24503.334 -> 21717
$ws.Cells.Item(135, 12).Value = 195453  # L135: 220530.006 -> 195453
$ws.Cells.Item(135, 14).Value = -200523  # N135: -225600.006 -> -200523

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(26, 8).Value = 34995  # H26: 0 -> 34995
$ws.Cells.Item(26, 10).Value = 34995  # J26: 0 -> 34995
$ws.Cells.Item(26, 12).Value = 34995  # L26: 0 -> 34995
$ws.Cells.Item(26, 14).Value = -35555  # N26: None -> -35555
$ws.Cells.Item(50, 8).Value = 34995  # H50: 0 -> 34995
$ws.Cells.Item(50, 10).Value = 34995  # J50: 0 -> 34995
$ws.Cells.Item(50, 12).Value = 34995  # L50: 0 -> 34995
$ws.Cells.Item(50, 14).Value = -35991  # N50: None -> -35991
$ws.Cells.Item(70, 8).Value = 5554  # H70: 5316.5713 -> 5554
$ws.Cells.Item(70, 9).Value = 5554  # I70: 5369.3335 -> 5554
$ws.Cells.Item(70, 10).Value = 0  # J70: 5000 -> 0
$ws.Cells.Item(70, 11).Value = 5554  # K70: 5369.3335 -> 5554
$ws.Cells.Item(70, 12).Value = 0  # L70: 5000 -> 0
$ws.Cells.Item(70, 13).Value = -5284  # M70: -5099.3335 -> -5284
$ws.Cells.Item(70, 14).ClearContents()  # N70: remove (was -5540)
$ws.Cells.Item(73, 8).Value = 5554  # H73: 5316.5713 -> 5554
$ws.Cells.Item(73, 9).Value = 5554  # I73: 5369.3335 -> 5554
$ws.Cells.Item(73, 10).Value = 0  # J73: 5000 -> 0
$ws.Cells.Item(73, 11).Value = 5554  # K73: 5369.3335 -> 5554
$ws.Cells.Item(73, 12).Value = 0  # L73: 5000 -> 0
$ws.Cells.Item(73, 13).Value = -4618  # M73: -4433.3335 -> -4618
$ws.Cells.Item(73, 14).ClearContents()  # N73: remove (was -6872)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 16711.572  # H7: 27758 -> 16711.572
$ws.Cells.Item(7, 9).Value = 3583  # I7: 3999.5 -> 3583
$ws.Cells.Item(7, 10).Value = 26558  # J7: 37261.4 -> 26558
$ws.Cells.Item(7, 11).Value = 3583  # K7: 3999.5 -> 3583
$ws.Cells.Item(7, 12).Value = 26558  # L7: 37261.4 -> 26558
$ws.Cells.Item(7, 13).Value = -3471  # M7: -3887.5 -> -3471
$ws.Cells.Item(7, 14).Value = -26782  # N7: -37485.4 -> -26782
$ws.Cells.Item(10, 8).Value = 115  # H10: 0 -> 115
$ws.Cells.Item(10, 9).Value = 30  # I10: 0 -> 30
$ws.Cells.Item(10, 10).Value = 200  # J10: 0 -> 200
$ws.Cells.Item(10, 11).Value = 30  # K10: 0 -> 30
$ws.Cells.Item(10, 12).Value = 200  # L10: 0 -> 200
$ws.Cells.Item(10, 13).Value = 110  # M10: None -> 110
$ws.Cells.Item(10, 14).Value = -480  # N10: None -> -480
$ws.Cells.Item(12, 8).Value = 5100  # H12: 10000 -> 5100
$ws.Cells.Item(12, 10).Value = 200  # J12: 0 -> 200
$ws.Cells.Item(12, 12).Value = 200  # L12: 0 -> 200
$ws.Cells.Item(12, 14).Value = -540  # N12: None -> -540
$ws.Cells.Item(54, 8).Value = 40084  # H54: 0 -> 40084
$ws.Cells.Item(54, 10).Value = 40084  # J54: 0 -> 40084
$ws.Cells.Item(54, 12).Value = 40084  # L54: 0 -> 40084
$ws.Cells.Item(54, 14).Value = -41372  # N54: None -> -41372
$ws.Cells.Item(100, 8).Value = 4671.4736  # H100: 4587.9 -> 4671.4736
$ws.Cells.Item(100, 9).Value = 3243.2856  # I100: 3227.0667 -> 3243.2856
$ws.Cells.Item(100, 11).Value = 3243.2856  # K100: 3227.0667 -> 3243.2856
$ws.Cells.Item(100, 13).Value = -2702.2856  # M100: -2686.0667 -> -2702.2856
$ws.Cells.Item(126, 8).Value = 16711.572  # H126: 27758 -> 16711.572
$ws.Cells.Item(126, 9).Value = 3583  # I126: 3999.5 -> 3583
$ws.Cells.Item(126, 10).Value = 26558  # J126: 37261.4 -> 26558
$ws.Cells.Item(126, 11).Value = 10749  # K126: 11998.5 -> 10749
$ws.Cells.Item(126, 12).Value = 79674  # L126: 111784.2 -> 79674
$ws.Cells.Item(126, 13).Value = -8279  # M126: -9528.5 -> -8279
$ws.Cells.Item(126, 14).Value = -84614  # N126: -116724.2 -> -84614
$ws.Cells.Item(136, 8).Value = 5103.32  # H136: 4838.643 -> 5103.32
$ws.Cells.Item(136, 9).Value = 3233.7646  # I136: 3143.65 -> 3233.7646
$ws.Cells.Item(136, 11).Value = 9701.293799999999  # K136: 9430.950000000001 -> 9701.293799999999
$ws.Cells.Item(136, 13).Value = -7151.293799999999  # M136: -6880.950000000001 -> -7151.293799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3397.4  # H81: 3554.3572 -> 3397.4
$ws.Cells.Item(81, 9).Value = 1996.8462  # I81: 2063.25 -> 1996.8462
$ws.Cells.Item(81, 11).Value = 3993.6924  # K81: 4126.5 -> 3993.6924
$ws.Cells.Item(81, 13).Value = -2932.6924  # M81: -3065.5 -> -2932.6924
$ws.Cells.Item(84, 8).Value = 3397.4  # H84: 3554.3572 -> 3397.4
$ws.Cells.Item(84, 9).Value = 1996.8462  # I84: 2063.25 -> 1996.8462
$ws.Cells.Item(84, 11).Value = 19968.462  # K84: 20632.5 -> 19968.462
$ws.Cells.Item(84, 13).Value = -14664.462  # M84: -15328.5 -> -14664.462
